$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.180.75'
$ws.Range("E2").Value = '  -0.64%  '

$ws.Range("D3").Value = '3.056.61'
$ws.Range("E3").Value = '  +1.22%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.441'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.21'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.81%  '

$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("E11").Value = '  +2.94%  '

$ws.Range("D12").Value = '3.579.19'
$ws.Range("E12").Value = '  +1.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.125'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000165'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.22%  '

$ws.Range("D16").Value = '57.138.99'
$ws.Range("E16").Value = '  -0.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.89%  '

$ws.Range("D18").Value = '3.057.63'
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '331.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.08%  '

$ws.Range("E22").Value = '  +0.49%  '

$ws.Range("E23").Value = '  +1.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.75%  '

$ws.Range("D25").Value = '3.177.23'
$ws.Range("E25").Value = '  +1.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.164'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.39%  '

$ws.Range("D28").Value = '0.0₃0890'
$ws.Range("E28").Value = '  -4.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.21%  '

$ws.Range("E31").Value = '  -0.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.98%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '151.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.19'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0671'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.14%  '

$ws.Range("B40").Value = 'RenzoRestakedETH'
$ws.Range("C40").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D40").Value = '3.095.27'
$ws.Range("E40").Value = '  +1.41%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.95%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.88'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.27%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.663'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.91%  '

$ws.Range("D45").Value = '2.201.67'
$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.963'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.62%  '

$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.40%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.61%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0241'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.59%  '

$ws.Range("E51").Value = '  +8.28%  '
